# Apply updated cryptos list values (prices + 1h volume deltas),
# matching the source diff cell-by-cell. A handful of D-column price
# strings are single-dot decimals (e.g. "20.30") that Excel's COM
# Value setter would otherwise auto-coerce to numbers (dropping the
# trailing zero); those are written with a leading apostrophe to force
# text, then the style is reset to "Normal" so no stray quote-prefix
# cell style lingers (keeps cells styleless, same as the source).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.122.16'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '1.624.15'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('D5').Value = "'214.11"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.50%  '
$ws.Range('E6').Value = '  +1.39%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -1.51%  '
$ws.Range('D9').Value = "'0.0627"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').Value = "'20.30"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('D12').Value = '1.622.76'
$ws.Range('E12').Value = '  -1.67%  '
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '27.103.17'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = "'64.60"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.30%  '
$ws.Range('D17').Value = '0.0₃0745'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = "'216.65"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.23%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').Value = "'6.93"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('E22').Value = '  -6.40%  '
$ws.Range('E23').Value = '  -1.88%  '
$ws.Range('D24').Value = "'148.09"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  -3.41%  '
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('D28').Value = "'15.61"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.25%  '
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('E31').Value = '  -0.77%  '
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('D33').Value = '1.343.69'
$ws.Range('E33').Value = '  +5.41%  '
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('E35').Value = '  -0.49%  '
$ws.Range('D36').Value = "'0.0177"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +1.42%  '
$ws.Range('D38').Value = "'0.857"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('D40').Value = "'0.804"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.69%  '
$ws.Range('D41').Value = "'65.59"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.99%  '
$ws.Range('E42').Value = '  -1.01%  '
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('D44').Value = '1.761.67'
$ws.Range('E44').Value = '  -1.21%  '
$ws.Range('E45').Value = '  -1.36%  '
$ws.Range('D46').Value = "'0.885"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +32.44%  '
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('E48').Value = '  -0.98%  '
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('E50').Value = '  +1.77%  '
$ws.Range('E51').Value = '  -1.10%  '
